$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.490588188171387
$ws.Range("B1").Value = 3.713589668273926
$ws.Range("C1").Value = 2.131559133529663
$ws.Range("D1").Value = 1.243484616279602
$ws.Range("E1").Value = 0.7562664151191711
